# Updates the "Price" (D) and "Volume(1h)" (E) columns for the crypto rows
# with freshly scraped values, mirroring the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.001.74'
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").Value = '2.350.43'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '''301.31'
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").Value = '''99.87'
$ws.Range("E6").Value = '  +2.54%  '
$ws.Range("E7").Value = '  -1.23%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '''0.514'
$ws.Range("E9").Value = '  -4.79%  '
$ws.Range("D10").Value = '''34.82'
$ws.Range("E10").Value = '  -3.17%  '
$ws.Range("D11").Value = '''0.0789'
$ws.Range("E11").Value = '  -2.98%  '
$ws.Range("E12").Value = '  -3.88%  '
$ws.Range("E13").Value = '  -1.40%  '
$ws.Range("D14").Value = '2.709.44'
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").Value = '2.351.73'
$ws.Range("E15").Value = '  +0.62%  '
$ws.Range("D16").Value = '''13.75'
$ws.Range("E16").Value = '  -3.14%  '
$ws.Range("D17").Value = '''0.809'
$ws.Range("E17").Value = '  -3.37%  '
$ws.Range("D18").Value = '45.992.20'
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("D19").Value = '''12.66'
$ws.Range("E19").Value = '  -7.56%  '
$ws.Range("D20").Value = '0.0₃0966'
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("E21").Value = '  -3.04%  '
$ws.Range("D22").Value = '''66.44'
$ws.Range("E22").Value = '  -2.04%  '
$ws.Range("D23").Value = '''245.25'
$ws.Range("E23").Value = '  -2.29%  '
$ws.Range("E24").Value = '  -5.62%  '
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("E26").Value = '  -4.63%  '
$ws.Range("D27").Value = '''40.91'
$ws.Range("E27").Value = '  -4.03%  '
$ws.Range("E28").Value = '  -3.79%  '
$ws.Range("D29").Value = '''9.71'
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("D30").Value = '''21.01'
$ws.Range("E30").Value = '  +3.35%  '
$ws.Range("D31").Value = '''3.65'
$ws.Range("E31").Value = '  +15.57%  '
$ws.Range("D32").Value = '''2.79'
$ws.Range("E32").Value = '  +6.70%  '
$ws.Range("E33").Value = '  -6.73%  '
$ws.Range("D34").Value = '''145.32'
$ws.Range("E34").Value = '  -1.14%  '
$ws.Range("D35").Value = '''0.0773'
$ws.Range("E35").Value = '  -5.48%  '
$ws.Range("E36").Value = '  -2.08%  '
$ws.Range("E37").Value = '  -2.95%  '
$ws.Range("D38").Value = '''1.80'
$ws.Range("E38").Value = '  -0.91%  '
$ws.Range("D39").Value = '''15.30'
$ws.Range("E39").Value = '  +8.79%  '
$ws.Range("D40").Value = '''3.87'
$ws.Range("E40").Value = '  -3.31%  '
$ws.Range("E41").Value = '  -4.62%  '
$ws.Range("E42").Value = '  -5.95%  '
$ws.Range("D43").Value = '''0.999'
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").Value = '1.850.03'
$ws.Range("E44").Value = '  +2.43%  '
$ws.Range("D45").Value = '''90.60'
$ws.Range("E45").Value = '  -1.07%  '
$ws.Range("D46").Value = '''1.82'
$ws.Range("E46").Value = '  -6.71%  '
$ws.Range("D47").Value = '''0.186'
$ws.Range("E47").Value = '  -4.88%  '
$ws.Range("D48").Value = '''70.99'
$ws.Range("E48").Value = '  -5.37%  '
$ws.Range("D49").Value = '2.580.74'
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").Value = '''8.01'
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").Value = '''96.05'
$ws.Range("E51").Value = '  -3.05%  '
